$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 534
$ws.Range("I2").Value = 1312
$ws.Range("J2").Value = 5696
$ws.Range("K2").Value = 30
$ws.Range("L2").Value = 1515
$ws.Range("M2").Value = 99
$ws.Range("N2").Value = 1099
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 22
$ws.Range("Q2").Value = 11
$ws.Range("R2").Value = 77
$ws.Range("S2").Value = 631
$ws.Range("T2").Value = 1041
$ws.Range("U2").Value = 71
$ws.Range("V2").Value = 8856
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 8782
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 137
$ws.Range("AA2").Value = 60
